# Update "mechanism" table values (stabilize cleaning files & internal
# validity analysis for number of pawns balance).
# All target cells live on the single active worksheet and hold cached
# string results of external-link formulas ("=[1]mechanism!..."); we only
# need to update the displayed/cached values, matching what Excel would
# show after the external workbook's source numbers changed and were
# refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C5"  = "9.76***"
    "C6"  = "(2.74)"
    "C7"  = "-0.58"
    "C8"  = "(2.23)"
    "C11" = "0.017"
    "C12" = "45.8"

    "B17" = "-0.071***"
    "C17" = "-4.12***"
    "H17" = "-0.20***"
    "C18" = "(1.28)"
    "H18" = "(0.050)"
    "B19" = "-0.027**"
    "C19" = "-1.82*"
    "D19" = "0.050*"
    "H19" = "-0.079*"
    "C20" = "(1.05)"
    "H20" = "(0.043)"
    "C22" = "2492"
    "D22" = "2492"
    "H22" = "2492"
    "B23" = "0.011"
    "C23" = "0.024"
    "D23" = "0.034"
    "H23" = "0.028"
    "C24" = "9.68"
    "D24" = "0.71"
}

# These cells hold cached *string* results of formulas that pull from an
# external workbook ("=[1]mechanism!..."). The external source isn't
# reachable here to do a live "update links" refresh, so we drive each
# cell with a literal-string formula (="9.76***") rather than a bare
# .Value assignment. That keeps the cell a formula/string cell (t="str",
# same cell style) instead of silently being re-typed into a number
# (dropping the parentheses/trailing zeros/asterisks formatting that carry
# significance in this regression-table text), which is what a plain
# .Value write would do.
foreach ($addr in $updates.Keys) {
    $text = $updates[$addr]
    $ws.Range($addr).Formula = '="' + $text + '"'
}
